# Fix duplicate shared-string entry: " Léodagan" (leading space) -> "Léodagan"
# This removes the stray leading space that was creating a near-duplicate
# shared string alongside the other entries (better modularisation & avoiding
# duplicates, per the commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "Léodagan"

# Reflect the new active cell/selection on the sheet (was E7, now B7)
$ws.Range("B7").Select()
